# Tidsplan_xjobb.xlsx update
# - Update PlanHours (P15) from 18 to 21
# - Update the pending-comment date in Q15 from "2020-01-14 Pending" to "2020-01-30 Pending"
# - Move active selection from Q18 to P15 (and scroll window accordingly)
# All dependent formulas (S15, P36, S36, P37, S37, P38, P43, P44, P45) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P15").Value = 21
$ws.Range("Q15").Value = "2020-01-30 Pending"

# Update the selected / active cell to match the new focus point of the edit (P15)
[void]$ws.Range("P15").Select()
